$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Update the "总计" (Total) sheet: insert a new 2022-Q4 summary row
#    before the existing 2022-Q3 / 2022-Q2 rows.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item(1)

# Shift the existing two data rows down by one row (row2->row3, row3->row4)
$total.Cells.Item(4,1).Value = 2
$total.Cells.Item(4,2).Value = "2022-Q2"
$total.Cells.Item(4,3).Value = 2
$total.Cells.Item(4,4).Value = 0.02

$total.Cells.Item(3,1).Value = 1
$total.Cells.Item(3,2).Value = "2022-Q3"
$total.Cells.Item(3,3).Value = 8
$total.Cells.Item(3,4).Value = 0.27

# Write the new 2022-Q4 summary row into row 2
$total.Cells.Item(2,1).Value = 0
$total.Cells.Item(2,2).Value = "2022-Q4"
$total.Cells.Item(2,3).Value = 3
$total.Cells.Item(2,4).Value = 0.98

# Row 4 is brand new territory on this sheet (dimension grows from D3 to D4);
# give its "A" cell the same style used by the other index cells (A2/A3).
$total.Range("A3").Copy()
$total.Range("A4").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 2. Create the new "2022-Q4" worksheet by duplicating the "2022-Q3"
#    sheet (so it keeps identical look & feel / styles) and placing it
#    right before "2022-Q3".
# ---------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")
$q3.Copy($q3)
$q4 = $wb.Worksheets.Item(2)
$q4.Name = "2022-Q4"

# Force the fund-code / numeric-looking text columns to stay text so
# leading zeros (e.g. "010363") and decimal strings are preserved exactly
# like the source data (which stores them as text, not numbers).
$q4.Range("B2:G4").NumberFormat = "@"

$q4.Cells.Item(2,2).Value = "010363"
$q4.Cells.Item(2,3).Value = "信澳匠心臻选两年持有期混合"
$q4.Cells.Item(2,4).Value = "37.44"
$q4.Cells.Item(2,5).Value = "92.71"
$q4.Cells.Item(2,6).Value = "2.46"
$q4.Cells.Item(2,7).Value = "0.9210"
$q4.Cells.Item(2,8).Value = 10

$q4.Cells.Item(3,2).Value = "013721"
$q4.Cells.Item(3,3).Value = "信澳景气优选混合A"
$q4.Cells.Item(3,4).Value = "1.08"
$q4.Cells.Item(3,5).Value = "92.37"
$q4.Cells.Item(3,6).Value = "3.89"
$q4.Cells.Item(3,7).Value = "0.0420"
$q4.Cells.Item(3,8).Value = 9

$q4.Cells.Item(4,2).Value = "013722"
$q4.Cells.Item(4,3).Value = "信澳景气优选混合C"
$q4.Cells.Item(4,4).Value = "0.48"
$q4.Cells.Item(4,5).Value = "92.37"
$q4.Cells.Item(4,6).Value = "3.89"
$q4.Cells.Item(4,7).Value = "0.0187"
$q4.Cells.Item(4,8).Value = 9

# The copied sheet has 9 rows of 2022-Q3 data (rows 2-9); the new
# 2022-Q4 sheet only needs 3 data rows (rows 2-4), so remove the rest.
$q4.Rows("5:9").Delete()
